# Apply crypto price/volume updates per the commit diff (Fri Nov 10 03:00:09 UTC 2023 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.873.74'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '2.123.47'
$ws.Range("E3").Value = '  +10.78%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.672'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.05%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.76'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '62.02'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.66%  '
$ws.Range("E10").Value = '  +2.90%  '
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").Value = '2.428.54'
$ws.Range("E13").Value = '  +10.56%  '
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.860'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.70%  '
$ws.Range("D16").Value = '2.122.30'
$ws.Range("E16").Value = '  +10.66%  '
$ws.Range("E17").Value = '  +1.56%  '
$ws.Range("D18").Value = '36.921.54'
$ws.Range("E18").Value = '  +3.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '242.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.97%  '
$ws.Range("E23").Value = '  +1.23%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '173.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.53'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.41%  '
$ws.Range("E29").Value = '  -8.55%  '
$ws.Range("E30").Value = '  -3.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +52.20%  '
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0966'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +20.72%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.40%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.921'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.89%  '
$ws.Range("E40").Value = '  -7.60%  '
$ws.Range("E41").Value = '  +8.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0225'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.69%  '
$ws.Range("E44").Value = '  +17.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.20%  '
$ws.Range("D46").Value = '1.369.43'
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0840'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.88%  '
$ws.Range("D48").Value = '2.324.14'
$ws.Range("E48").Value = '  +10.95%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.97'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.46%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("E51").Value = '  +2.13%  '
